$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spec")

$ws.Range("F2").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 5
$ws.Range("F7").Value = 10
$ws.Range("F8").Value = 10
$ws.Range("F9").Value = 10
$ws.Range("F10").Value = 10
$ws.Range("F11").Value = 10
$ws.Range("F12").Value = 10
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = 10
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 10
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 10
$ws.Range("F20").Value = 10
$ws.Range("F21").Value = 10
$ws.Range("F22").Value = 10
$ws.Range("F23").Value = 10
$ws.Range("F24").Value = 15
$ws.Range("F26").Value = 20
$ws.Range("F27").Value = 5
$ws.Range("F29").Value = 10
$ws.Range("F32").Value = 5
$ws.Range("F33").Value = 10
$ws.Range("F34").Value = 5
$ws.Range("F38").Value = 5
$ws.Range("F39").Value = 5
$ws.Range("F40").Value = 5
$ws.Range("F42").Value = 5
$ws.Range("F43").Value = 5
$ws.Range("F44").Value = 10
$ws.Range("F45").Value = 5
$ws.Range("F46").Value = 5
$ws.Range("F47").Value = 15
$ws.Range("F48").Value = 5
$ws.Range("F51").Value = 10
$ws.Range("F52").Value = 5
$ws.Range("F53").Value = 10
$ws.Range("F54").Value = 5
$ws.Range("F55").Value = 10
$ws.Range("F56").Value = 5
$ws.Range("F57").Value = 5
$ws.Range("F60").Value = 10
$ws.Range("F61").Value = 10
$ws.Range("F62").Value = 15
$ws.Range("F63").Value = 15
$ws.Range("F64").Value = 15
$ws.Range("F65").Value = 15
$ws.Range("F66").Value = 15
$ws.Range("F70").Value = 5
$ws.Range("F71").Value = 10
$ws.Range("F73").Value = 10
$ws.Range("F74").Value = 5
$ws.Range("F76").Value = 20

# Delete row 77 entirely (shift cells up), removing the now-obsolete "Soundprof" entry
$ws.Rows.Item(77).Delete()
